# Insert a new data row at row 36 in the single worksheet of the workbook.
# This shifts the existing rows 36-100 down to rows 37-101 (and updates the
# sheet dimension automatically), then populates the newly inserted row 36
# with its own data (mirrors the row immediately below it except for the
# columns that differ: Fecha, Volumen, Precio minimo/maximo/promedio,
# Unidad de comercializacion and Precio $/Kg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 36 (old row 36 -> new row 37, etc.)
$ws.Rows("36:36").Insert()

# Populate the new row 36 with its data.
$ws.Cells.Item(36, 1).Value  = 5
$ws.Cells.Item(36, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(36, 3).Value  = "Maule"
$ws.Cells.Item(36, 4).Value  = 44498
$ws.Cells.Item(36, 5).Value  = 7
$ws.Cells.Item(36, 6).Value  = 100112031
$ws.Cells.Item(36, 7).Value  = "Poroto verde"
$ws.Cells.Item(36, 8).Value  = "Sin especificar"
$ws.Cells.Item(36, 9).Value  = "Primera"
$ws.Cells.Item(36, 10).Value = 100
$ws.Cells.Item(36, 11).Value = 50000
$ws.Cells.Item(36, 12).Value = 50000
$ws.Cells.Item(36, 13).Value = 50000
$ws.Cells.Item(36, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(36, 15).Value = "Región del Maule"
$ws.Cells.Item(36, 16).Value = 2000
$ws.Cells.Item(36, 17).Value = 25
$ws.Cells.Item(36, 18).Value = "Hortaliza"
